# Insert a new data row at row 30 (weekly Mango price record for
# Macroferia Regional de Talca), shifting all subsequent rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Inserting the whole row pushes existing rows 30..109 down to 31..110
# and extends the used range to A1:T110, exactly like Excel's
# Rows("30").Insert() behavior (shift-down, inherits formatting from the
# row above).
$ws.Rows("30").Insert()

# Populate the newly inserted row 30 with the new record's data.
$ws.Range("A30").Value = 5
$ws.Range("B30").Value = "Macroferia Regional de Talca"
$ws.Range("C30").Value = "Maule"
$ws.Range("D30").Value = 44607
$ws.Range("E30").Value = 7
$ws.Range("F30").Value = "Fruta"
$ws.Range("G30").Value = 100108
$ws.Range("H30").Value = "Tropicales y subtropicales"
$ws.Range("I30").Value = 100108002
$ws.Range("J30").Value = "Mango"
$ws.Range("K30").Value = "Sin especificar"
$ws.Range("L30").Value = "Primera"
$ws.Range("M30").Value = 240
$ws.Range("N30").Value = 7000
$ws.Range("O30").Value = 7000
$ws.Range("P30").Value = 7000
$ws.Range("Q30").Value = "$/bandeja 4 kilos"
$ws.Range("R30").Value = "Perú"
$ws.Range("S30").Value = 1750
$ws.Range("T30").Value = 4
